$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range("E2").Value = "2026-02-10 06:18:55"
$ws.Range("E3").Value = "2026-02-10 06:18:57"
$ws.Range("G3").Value = "184 cm"
$ws.Range("I3").Value = "7.7 mm"
$ws.Range("E4").Value = "2026-02-10 06:18:59"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "86%"
$ws.Range("E5").Value = "2026-02-10 06:19:02"
$ws.Range("I5").Value = "12.3 mm"
$ws.Range("E6").Value = "2026-02-10 06:19:04"
$ws.Range("O6").Value = "7.3 °C"
$ws.Range("E7").Value = "2026-02-10 06:19:06"
$ws.Range("E8").Value = "2026-02-10 06:19:09"
$ws.Range("N8").Value = "8.2 °C 5:57 TU"
$ws.Range("E9").Value = "2026-02-10 06:19:11"
$ws.Range("O9").Value = "6.5 °C"
$ws.Range("E10").Value = "2026-02-10 06:19:14"
$ws.Range("N10").Value = "4.6 °C 5:59 TU"
$ws.Range("O10").Value = "7.0 °C"
$ws.Range("E11").Value = "2026-02-10 06:19:16"
$ws.Range("O11").Value = "3.0 °C"
$ws.Range("E12").Value = "2026-02-10 06:19:19"
$ws.Range("N12").Value = "5.3 °C 5:31 TU"
$ws.Range("O12").Value = "6.8 °C"
$ws.Range("E13").Value = "2026-02-10 06:19:21"
$ws.Range("I13").Value = "2.1 mm"
$ws.Range("J13").Value = "1008.2 hPa"
$ws.Range("N13").Value = "2.4 °C 5:51 TU"
$ws.Range("E14").Value = "2026-02-10 06:19:23"
$ws.Range("O14").Value = "9.4 °C"
$ws.Range("E15").Value = "2026-02-10 06:19:26"
$ws.Range("O15").Value = "6.3 °C"
$ws.Range("E16").Value = "2026-02-10 06:19:28"
$ws.Range("I16").Value = "12.2 mm"
$ws.Range("E17").Value = "2026-02-10 06:19:31"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "91%"
$ws.Range("O17").Value = "2.8 °C"
$ws.Range("E18").Value = "2026-02-10 06:19:33"
$ws.Range("N18").Value = "4.7 °C 5:57 TU"
$ws.Range("O18").Value = "7.3 °C"
$ws.Range("E19").Value = "2026-02-10 06:19:36"
$ws.Range("N19").Value = "3.2 °C 5:54 TU"
$ws.Range("E20").Value = "2026-02-10 06:19:39"
$ws.Range("G20").Value = "123 cm"
$ws.Range("O20").Value = "-1.1 °C"
$ws.Range("E21").Value = "2026-02-10 06:19:41"
$ws.Range("I21").Value = "3.2 mm"
$ws.Range("E22").Value = "2026-02-10 06:19:43"
$ws.Range("I22").Value = "0.3 mm"
$ws.Range("M22").Value = "-1.4 °C 5:48 TU"
$ws.Range("O22").Value = "-1.9 °C"
$ws.Range("E23").Value = "2026-02-10 06:19:46"
$ws.Range("I23").Value = "10.6 mm"
$ws.Range("E24").Value = "2026-02-10 06:19:48"
$ws.Range("I24").Value = "1.5 mm"
$ws.Range("J24").Value = "1007.0 hPa"
$ws.Range("N24").Value = "8.1 °C 5:56 TU"
$ws.Range("E25").Value = "2026-02-10 06:19:50"
$ws.Range("G25").Value = "119 cm"
$ws.Range("I25").Value = "6.7 mm"
$ws.Range("E26").Value = "2026-02-10 06:19:53"
$ws.Range("M26").Value = "4.6 °C 5:53 TU"
$ws.Range("O26").Value = "3.0 °C"
$ws.Range("E27").Value = "2026-02-10 06:19:55"
$ws.Range("I27").Value = "1.9 mm"
$ws.Range("L27").Value = "32.8 km/h - 246º 5:55 TU"
$ws.Range("E28").Value = "2026-02-10 06:19:57"
$ws.Range("J28").Value = "1005.6 hPa"
$ws.Range("N28").Value = "3.3 °C 5:59 TU"
$ws.Range("O28").Value = "5.0 °C"
$ws.Range("E29").Value = "2026-02-10 06:20:00"
$ws.Range("H29").Value = ""
$ws.Range("I29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = ""
$ws.Range("O29").Value = ""
$ws.Range("E30").Value = "2026-02-10 06:20:16"
$ws.Range("N30").Value = "6.5 °C 5:55 TU"
$ws.Range("E31").Value = "2026-02-10 06:20:32"
$ws.Range("J31").Value = "1004.5 hPa"
$ws.Range("E32").Value = "2026-02-10 06:20:35"
$ws.Range("L32").Value = "30.6 km/h - 295º 5:51 TU"
$ws.Range("O32").Value = "7.7 °C"
$ws.Range("E33").Value = "2026-02-10 06:20:37"
$ws.Range("I33").Value = "4.7 mm"
$ws.Range("N33").Value = "1.7 °C 5:32 TU"
$ws.Range("E34").Value = "2026-02-10 06:20:40"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "80%"
$ws.Range("I34").Value = "2.2 mm"
$ws.Range("N34").Value = "0.9 °C 5:58 TU"
$ws.Range("O34").Value = "2.5 °C"
$ws.Range("E35").Value = "2026-02-10 06:20:43"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "83%"
$ws.Range("J35").Value = "1005.6 hPa"
$ws.Range("E36").Value = "2026-02-10 06:20:46"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "94%"
$ws.Range("J36").Value = "1005.0 hPa"
$ws.Range("N36").Value = "6.7 °C 5:46 TU"
$ws.Range("O36").Value = "8.8 °C"
$ws.Range("E37").Value = "2026-02-10 06:20:49"
$ws.Range("E38").Value = "2026-02-10 06:20:52"
$ws.Range("N38").Value = "6.5 °C 5:56 TU"
$ws.Range("O38").Value = "7.6 °C"
$ws.Range("E39").Value = "2026-02-10 06:20:54"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "86%"
$ws.Range("I39").Value = "2.6 mm"
$ws.Range("O39").Value = "-0.2 °C"
$ws.Range("E40").Value = "2026-02-10 06:20:57"
$ws.Range("I40").Value = "3.8 mm"
$ws.Range("N40").Value = "4.1 °C 5:36 TU"
$ws.Range("E41").Value = "2026-02-10 06:21:00"
$ws.Range("E42").Value = "2026-02-10 06:21:03"
$ws.Range("N42").Value = "6.9 °C 5:52 TU"
$ws.Range("O42").Value = "8.0 °C"
$ws.Range("E43").Value = "2026-02-10 06:21:05"
$ws.Range("E44").Value = "2026-02-10 06:21:07"
$ws.Range("I44").Value = "7.3 mm"
$ws.Range("O44").Value = "-0.3 °C"
$ws.Range("E45").Value = "2026-02-10 06:21:10"
$ws.Range("I45").Value = "15.9 mm"
$ws.Range("M45").Value = "4.2 °C 5:59 TU"
$ws.Range("E46").Value = "2026-02-10 06:21:13"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "100%"
$ws.Range("J46").Value = "1006.8 hPa"
$ws.Range("L46").Value = "12.6 km/h - 51º 5:30 TU"
